$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D data cells (D2:D51) to text format so that
# numeric-looking price strings (e.g. "0.9965") are not auto-converted to
# numbers by Excel's type inference, matching the original inline-string cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.988.15"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.829.65"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").Value = "0.9965"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "243.76"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "0.6321"
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("D7").Value = "0.9985"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "0.07523"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "0.2943"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "22.97"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "0.07726"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "1.826.68"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "4.995"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "0.6713"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "83.14"
$ws.Range("E15").Value = "  +1.05%  "
$ws.Range("D16").Value = "0.000009721"
$ws.Range("E16").Value = "  +7.40%  "
$ws.Range("D17").Value = "6.078"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "29.039.54"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "12.56"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "226.54"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "0.9978"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "7.182"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "0.9979"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "159.78"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "0.1403"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").Value = "8.553"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").Value = "17.93"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "1.497"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "4.122"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "4.082"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "1.202"
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "0.05379"
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("D33").Value = "1.864"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").Value = "0.7442"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "2.655"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").Value = "1.243.81"
$ws.Range("E37").Value = "  -2.29%  "
$ws.Range("D38").Value = "0.01787"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").Value = "2.750"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "6.598"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").Value = "0.9053"
$ws.Range("E41").Value = "  +1.48%  "
$ws.Range("D42").Value = "0.9981"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "101.70"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "1.983.21"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  +3.68%  "
$ws.Range("D46").Value = "64.88"
$ws.Range("E46").Value = "  +2.21%  "
$ws.Range("D47").Value = "0.5097"
$ws.Range("E47").Value = "  -0.34%  "
$ws.Range("D48").Value = "0.4078"
$ws.Range("E48").Value = "  +2.82%  "
$ws.Range("D49").Value = "9.069"
$ws.Range("E49").Value = "  +2.57%  "
$ws.Range("D50").Value = "6.770"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("D51").Value = "0.05765"
$ws.Range("E51").Value = "  +0.16%  "

# Restore the default (unstyled) cell style now that the text values are
# safely stored, so the cells match the original workbook's formatting.
$ws.Range("D2:D51").Style = "Normal"
